$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ "A"="1354D9E0-5D89-44D6-A4BB-E93417D95B5F"; "B"="P-6585"; "C"="2025-08-22T13:26:42"; "D"="ADRIANA QUIROZ RODRIGUEZ"; "E"="QURA720718KP7"; "F"="010 000 0246 00 00 PROPOFOL. EMULSIÓN INYECTABLE.EMULSION INYECTABLE CADA AMPOLLETA O FRASCO AMPULA CONTIENE: EMULSION CON O SIN EDETATO DISODICO (DIHIDRATADO) 200 MG. ENVASE CON 5 AMPOLLETAS O FRASCOS AMPULA DE 20 ML."; "G"=250; "H"=23750; "I"="P-6585 ADRIANA QUIROZ RODRIGUEZ.xml" },
    @{ "A"="F27BFC72-8E1B-404C-BAAF-D950102B1B09"; "B"="P-6570"; "C"="2025-08-22T11:15:31"; "D"="INSTITUTO MEXICANO DEL SEGURO SOCIAL"; "E"="IMS421231I45"; "F"="010 000 4411 00 00 LATANOPROSTSOLUCION OFTALMICA CADA ML CONTIENE: LATANOPROST 50 MICROGRAMOS ENVASE CON UN FRASCO GOTERO CON 2.5 ML."; "G"=1996; "H"=165668; "I"="P-6570 IMSS.xml" },
    @{ "A"="69D4D535-CF02-44A5-9A3D-2A632469B947"; "B"="P-6571"; "C"="2025-08-22T11:23:32"; "D"="INSTITUTO MEXICANO DEL SEGURO SOCIAL"; "E"="IMS421231I45"; "F"="010 000 4411 00 00 LATANOPROSTSOLUCION OFTALMICA CADA ML CONTIENE: LATANOPROST 50 MICROGRAMOS ENVASE CON UN FRASCO GOTERO CON 2.5 ML."; "G"=2146; "H"=178118; "I"="P-6571 IMSS.xml" },
    @{ "A"="9141F3EB-728B-47B4-BF84-FCEBBAC8BBE0"; "B"="P-6572"; "C"="2025-08-22T11:28:30"; "D"="INSTITUTO MEXICANO DEL SEGURO SOCIAL"; "E"="IMS421231I45"; "F"="010 000 4411 00 00 LATANOPROSTSOLUCION OFTALMICA CADA ML CONTIENE: LATANOPROST 50 MICROGRAMOS ENVASE CON UN FRASCO GOTERO CON 2.5 ML."; "G"=10; "H"=830; "I"="P-6572 IMSS.xml" },
    @{ "A"="43BC1517-6DA3-42FC-A35E-427E974F782B"; "B"="P-6573"; "C"="2025-08-22T11:34:32"; "D"="INSTITUTO MEXICANO DEL SEGURO SOCIAL"; "E"="IMS421231I45"; "F"="010 000 4411 00 00 LATANOPROSTSOLUCION OFTALMICA CADA ML CONTIENE: LATANOPROST 50 MICROGRAMOS ENVASE CON UN FRASCO GOTERO CON 2.5 ML."; "G"=83; "H"=6889; "I"="P-6573 IMSS.xml" },
    @{ "A"="0658F2A5-56B5-441C-B17B-D3246DF54929"; "B"="P-6574"; "C"="2025-08-22T11:45:14"; "D"="INSTITUTO MEXICANO DEL SEGURO SOCIAL"; "E"="IMS421231I45"; "F"="010 000 4411 00 00 LATANOPROSTSOLUCION OFTALMICA CADA ML CONTIENE: LATANOPROST 50 MICROGRAMOS ENVASE CON UN FRASCO GOTERO CON 2.5 ML."; "G"=25; "H"=2075; "I"="P-6574 IMSS.xml" },
    @{ "A"="C20465C7-5DAC-48D2-AC0F-EE780230D57A"; "B"="P-6575"; "C"="2025-08-22T11:51:56"; "D"="INSTITUTO MEXICANO DEL SEGURO SOCIAL"; "E"="IMS421231I45"; "F"="010 000 4411 00 00 LATANOPROSTSOLUCION OFTALMICA CADA ML CONTIENE: LATANOPROST 50 MICROGRAMOS ENVASE CON UN FRASCO GOTERO CON 2.5 ML."; "G"=4589; "H"=380887; "I"="P-6575 IMSS.xml" },
    @{ "A"="EE275B25-6ED4-4E22-86D5-FAB968BE5FDA"; "B"="P-6576"; "C"="2025-08-22T11:58:58"; "D"="INSTITUTO MEXICANO DEL SEGURO SOCIAL"; "E"="IMS421231I45"; "F"="010 000 4411 00 00 LATANOPROSTSOLUCION OFTALMICA CADA ML CONTIENE: LATANOPROST 50 MICROGRAMOS ENVASE CON UN FRASCO GOTERO CON 2.5 ML."; "G"=1124; "H"=93292; "I"="P-6576 IMSS.xml" },
    @{ "A"="870D6D18-E68B-4820-90D8-B10DCFF903A1"; "B"="P-6577"; "C"="2025-08-22T12:13:55"; "D"="INSTITUTO MEXICANO DEL SEGURO SOCIAL"; "E"="IMS421231I45"; "F"="010 000 4411 00 00 LATANOPROSTSOLUCION OFTALMICA CADA ML CONTIENE: LATANOPROST 50 MICROGRAMOS ENVASE CON UN FRASCO GOTERO CON 2.5 ML."; "G"=78; "H"=6474; "I"="P-6577 IMSS.xml" },
    @{ "A"="5E63D996-999D-46C3-AFD2-DDCE228778D9"; "B"="P-6578"; "C"="2025-08-22T12:22:18"; "D"="INSTITUTO MEXICANO DEL SEGURO SOCIAL"; "E"="IMS421231I45"; "F"="010 000 4411 00 00 LATANOPROSTSOLUCION OFTALMICA CADA ML CONTIENE: LATANOPROST 50 MICROGRAMOS ENVASE CON UN FRASCO GOTERO CON 2.5 ML."; "G"=1466; "H"=121678; "I"="P-6578 IMSS.xml" },
    @{ "A"="3159F584-4DF8-4738-AC55-E24CF40439E8"; "B"="P-6579"; "C"="2025-08-22T12:28:30"; "D"="INSTITUTO MEXICANO DEL SEGURO SOCIAL"; "E"="IMS421231I45"; "F"="010 000 4411 00 00 LATANOPROSTSOLUCION OFTALMICA CADA ML CONTIENE: LATANOPROST 50 MICROGRAMOS ENVASE CON UN FRASCO GOTERO CON 2.5 ML."; "G"=997; "H"=82751; "I"="P-6579 IMSS.xml" },
    @{ "A"="0F825E1A-03C5-40D7-80FD-938FC390CB0E"; "B"="P-6580"; "C"="2025-08-22T12:43:44"; "D"="INSTITUTO MEXICANO DEL SEGURO SOCIAL"; "E"="IMS421231I45"; "F"="010 000 4411 00 00 LATANOPROSTSOLUCION OFTALMICA CADA ML CONTIENE: LATANOPROST 50 MICROGRAMOS ENVASE CON UN FRASCO GOTERO CON 2.5 ML."; "G"=7300; "H"=605900; "I"="P-6580 IMSS.xml" },
    @{ "A"="5D8CB6F1-C222-46D1-84F8-2C2AF908F49B"; "B"="P-6581"; "C"="2025-08-22T12:48:23"; "D"="INSTITUTO MEXICANO DEL SEGURO SOCIAL"; "E"="IMS421231I45"; "F"="010 000 4411 00 00 LATANOPROSTSOLUCION OFTALMICA CADA ML CONTIENE: LATANOPROST 50 MICROGRAMOS ENVASE CON UN FRASCO GOTERO CON 2.5 ML."; "G"=7031; "H"=583573; "I"="P-6581 IMSS.xml" },
    @{ "A"="6D3D9676-C144-4E4E-A37C-F7E34ECAA50B"; "B"="P-6582"; "C"="2025-08-22T12:52:50"; "D"="INSTITUTO MEXICANO DEL SEGURO SOCIAL"; "E"="IMS421231I45"; "F"="010 000 4411 00 00 LATANOPROSTSOLUCION OFTALMICA CADA ML CONTIENE: LATANOPROST 50 MICROGRAMOS ENVASE CON UN FRASCO GOTERO CON 2.5 ML."; "G"=5067; "H"=420561; "I"="P-6582 IMSS.xml" },
    @{ "A"="4A28D11B-AD4A-401C-AE95-40A08BD9F166"; "B"="P-6583"; "C"="2025-08-22T12:57:42"; "D"="INSTITUTO MEXICANO DEL SEGURO SOCIAL"; "E"="IMS421231I45"; "F"="010 000 4411 00 00 LATANOPROSTSOLUCION OFTALMICA CADA ML CONTIENE: LATANOPROST 50 MICROGRAMOS ENVASE CON UN FRASCO GOTERO CON 2.5 ML."; "G"=1452; "H"=120516; "I"="P-6583 IMSS.xml" },
    @{ "A"="B0C51690-3829-4651-B664-769AB6FD94CC"; "B"="P-6584"; "C"="2025-08-22T13:01:11"; "D"="INSTITUTO MEXICANO DEL SEGURO SOCIAL"; "E"="IMS421231I45"; "F"="010 000 4411 00 00 LATANOPROSTSOLUCION OFTALMICA CADA ML CONTIENE: LATANOPROST 50 MICROGRAMOS ENVASE CON UN FRASCO GOTERO CON 2.5 ML."; "G"=713; "H"=59179; "I"="P-6584 IMSS.xml" }
)

$startRow = 182
$cols = @("A","B","C","D","E","F","G","H","I")
$numericCols = @("G","H")

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    foreach ($col in $cols) {
        $addr = "$col$r"
        if ($numericCols -contains $col) {
            $ws.Range($addr).Value = [double]$rowData[$col]
        } else {
            $ws.Range($addr).Value = [string]$rowData[$col]
        }
    }
}

Write-Host "Done. UsedRange rows:"
Write-Host $ws.UsedRange.Rows.Count